# Update the Questions worksheet: replace the legacy Python-repr dump that
# lived in the shared string referenced from A2 with the pretty-printed JSON
# version, write it into A1 (dropping A1's old "0" placeholder value and its
# bold/bordered header style), then remove the now-unused second row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = 'questions = [
    {
        "title": "You have a database with employee ID numbers and you are trying to locate a person with a specific ID number. To access the location of an ID takes a constant amount of time. What is the time complexity of locating the ID?",
        "ques_type": 2,
        "options": [
            "O(1)",
            "O(0)",
            "\u03a9(0)",
            "\u03a9(c)"
        ],
        "score": "O(1)"
    },
    {
        "title": "You are trying to draw the recursion tree of an algorithm. The algorithm does a recursive call to the left or the right neighbor of the current number called from the sequence at each step.  What is the greatest number of branches each node in the recursive tree can have?",
        "ques_type": 2,
        "options": [
            "1",
            "2",
            "3",
            "4"
        ],
        "score": "2"
    },
    {
        "title": "You are trying to pick people to form a team for a project. You want to minimize the overall salary requirements of the team. There is only a limit on how many people can be picked from each department. You are using a greedy algorithm. What is the greedy choice property for picking the members of the team?",
        "ques_type": 2,
        "options": [
            "Pick the people with the lowest salary available from randomly picked departments.",
            "Pick the people with the lowest salary available from the department having the lowest average salary.",
            "Pick the people with the lowest salary available from the department having the lowest median salary.",
            "Pick the person with the lowest salary available from each of the departments."
        ],
        "score": "Pick the person with the lowest salary available from each of the departments."
    },
    {
        "title": "You are applying breadth-first search on a connected graph, where vertices represent different departments in the company and edges represent scheduled meetings between those departments. The graph has V vertices and E edges. What is the time complexity of listing all the meetings by traversing the graph?",
        "ques_type": 2,
        "options": [
            "O(V)",
            "O(E)",
            "O(V+E)",
            "O(V*E)"
        ],
        "score": "O(V+E)"
    }
]'

# Strip A1's existing bold + thin-border + center/top-aligned style before
# writing the new value (clearing formats *after* the value write perturbs
# the font size on this host, so the order here matters).
$ws.Range("A1").ClearFormats()
$ws.Range("A1").Value = $newText

# A2 (the old shared-string cell) is no longer needed once A1 carries the
# updated text, so drop the whole row and let the sheet shrink to A1 only.
$ws.Rows("2:2").Delete()
